{"js": "// Update the single-column results table in place: several summary\n// cells get new numbers, and the three \"raw sample\" rows (which held\n// a long tab-separated run of per-iteration numbers) are collapsed\n// down to the single headline number that used to live at the top of\n// the table.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// 0-based (rowIndex, columnIndex) -> new cell text.\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"312\",\n  4: \"0.00002\",\n  5: \"0.00069\",\n  6: \"0.00017\",\n  8: \"0.00031\",\n  9: \"0.00036\",\n  10: \"0.00040\",\n  11: \"0.06490\",\n  43: \"99.65\",\n  44: \"0.06\",\n  45: \"18\",\n};\n\nfor (const [rowIndex, newText] of Object.entries(updates)) {\n  const cell = table.getCell(Number(rowIndex), 0);\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Update the single-column results table in place: several summary\n# cells get new numbers, and the three \"raw sample\" rows (which held\n# a long tab-separated run of per-iteration numbers) are collapsed\n# down to the single headline number that used to live at the top of\n# the table.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# 1-based (row, column) -> new cell text (Word COM table indices are 1-based).\n$updates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"312\"\n    5  = \"0.00002\"\n    6  = \"0.00069\"\n    7  = \"0.00017\"\n    9  = \"0.00031\"\n    10 = \"0.00036\"\n    11 = \"0.00040\"\n    12 = \"0.06490\"\n    44 = \"99.65\"\n    45 = \"0.06\"\n    46 = \"18\"\n}\n\nforeach ($row in $updates.Keys) {\n    $cell = $t.Cell($row, 1)\n    $cell.Range.Text = $updates[$row]\n}\n"}
